# Tweak analytics opt out prompt:
# - Swap order so the "disable analytics" (opt out) option comes first,
#   followed by "enable analytics".
# - Reword the button captions and the prompt text itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A104").Value = "disable analytics"
$ws.Range("B104").Value = "Opt out"
$ws.Range("A105").Value = "disable analytics gamepad"
$ws.Range("B105").Value = "[B] Opt out"
$ws.Range("A106").Value = "enable analytics"
$ws.Range("B106").Value = "It's fine"
$ws.Range("A107").Value = "enable analytics gamepad"
$ws.Range("B107").Value = "[A] It's fine"
$ws.Range("A108").Value = "analytics prompt"
$ws.Range("B108").Value = "Opt out of anonymous crash reports and analytics?"

# Restore the view state to match where this edit was made: the sheet
# scrolled down and the "enable analytics" key cell selected.
$ws.Range("A106").Select()
